# Fruta / hortaliza, semanal
# Rotate the weekly data blocks (row pairs) down the sheet:
#   rows 4-5   <- old rows 10-11
#   rows 6-7   <- old rows 14-15
#   rows 8-9   <- old rows 12-13
#   rows 10-11 <- old rows 4-5
#   rows 12-13 <- old rows 6-7
#   rows 14-15 <- old rows 8-9
# Only columns D (Fecha), M (Volumen), N/O/P (Precios), R (Origen) and S (Precio $/Kg)
# actually change values; capture them first, then write back in the new order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "M", "N", "O", "P", "R", "S")

# Snapshot current ("before") values for the affected rows.
$snapshot = @{}
foreach ($row in 4..15) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowData
}

# Destination row -> source row mapping (the rotation observed in the diff).
$mapping = @{
    4  = 10
    5  = 11
    6  = 14
    7  = 15
    8  = 12
    9  = 13
    10 = 4
    11 = 5
    12 = 6
    13 = 7
    14 = 8
    15 = 9
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcData = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value = $srcData[$col]
    }
}
